$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price and Volume columns remain plain text so values such as
# "26.006.47" or "0.00000000113" are not reinterpreted as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.006.47"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "1.666.38"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "216.44"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").Value = "0.5096"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "0.2637"
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").Value = "0.06381"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").Value = "21.92"
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").Value = "0.07413"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "1.665.35"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").Value = "4.498"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "0.5825"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "0.000008470"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "64.12"
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("D17").Value = "26.053.60"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").Value = "4.924"
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "10.71"
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").Value = "189.23"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("D22").Value = "6.205"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "144.92"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "7.596"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").Value = "0.1191"
$ws.Range("E26").Value = "  +3.40%  "
$ws.Range("D27").Value = "15.62"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "0.06610"
$ws.Range("E28").Value = "  +15.79%  "
$ws.Range("D30").Value = "1.312"
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "3.504"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "1.630"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").Value = "1.015"
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("D35").Value = "0.6072"
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("D36").Value = "2.366"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").Value = "6.195"
$ws.Range("E38").Value = "  +5.70%  "
$ws.Range("D39").Value = "0.01605"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "1.074.15"
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("D41").Value = "0.8580"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").Value = "100.41"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").Value = "1.811.80"
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").Value = "0.00000000113"
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("D46").Value = "56.21"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "1.007"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "8.010"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").Value = "0.05206"
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").Value = "0.4287"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "5.949"
$ws.Range("E51").Value = "  +2.62%  "
